# "added Datadriven for test classes"
# Populate Sheet1 with the test-data used by the data-driven test classes
# (email subscription test + product search test), turning the two sample
# e-mail addresses into mailto: hyperlinks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- test_email_is_subscribed -------------------------------------------------
$ws.Range("A1").Value = "test_email_is_subscribed"
$ws.Range("A2").Value = "random_email"

$ws.Hyperlinks.Add($ws.Range("A3"), "mailto:dfsdashdgh@gh.lo", "", "", "dfsdashdgh@gh.lo") | Out-Null
$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:jhsjdsu@tyh.mk", "", "", "jhsjdsu@tyh.mk") | Out-Null

# --- test_search_product -------------------------------------------------------
$ws.Range("A6").Value = "test_search_product"
$ws.Range("A7").Value = "product_name"
$ws.Range("A8").Value = "top"

# Leave the selection where it was left in the authored workbook.
$ws.Range("F7").Select() | Out-Null
